# Generate Report for Handback
# - Update the "Status" text (shared across Overview + per-language sheets) to
#   reflect that the handback is now in sync with en-US.
# - Refresh the "Latest Handback DateTime" timestamps for zh-cn / de-de.
# - Clear the stale "Error Detail" text now that the handback files are current.
# - Widen the datetime-ish columns (C on the language sheets / E+F on Overview)
#   and shrink the now-mostly-empty "Error Detail" column (P) on the language
#   sheets to better fit their content.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column updates (shared string used in all four cells) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Latest Handback DateTime updates ---
$zhcn.Range("K2").Value = "2016-09-05 00:56:24"
$dede.Range("K2").Value = "2016-09-05 00:56:31"

# --- Error Detail is no longer applicable now the handback is current ---
$zhcn.Range("P2").Value = ""
$dede.Range("P2").Value = ""

# --- Column width adjustments ---
# NOTE: Excel's ColumnWidth setter snaps to whole-pixel increments (i.e. the
# stored "character width" only takes values of the form n/6 + 5/6), so the
# assigned values below are chosen to land on the closest attainable width to
# the desired 29.9777047293527 / 13.7470528738839 "character" widths.
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366   # E -> ~29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366   # F -> ~29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.144371396019366        # C -> ~29.9777047293527
$zhcn.Columns.Item(16).ColumnWidth = 12.913719540550566       # P -> ~13.7470528738839

$dede.Columns.Item(3).ColumnWidth = 29.144371396019366        # C -> ~29.9777047293527
$dede.Columns.Item(16).ColumnWidth = 12.913719540550566       # P -> ~13.7470528738839
